$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: NATHAN / DANSKIN / 2222 / ADMIN  ->  ALANA / TEST / 3333 / STAFF
$ws.Range("A3").Value = "ALANA"
$ws.Range("B3").Value = "TEST"
$ws.Range("C3").NumberFormat = "@"
$ws.Range("C3").Value = "3333"
$ws.Range("D3").Value = "STAFF"

# Row 4: ALANA / DANSKIN / 3333 / MANAGER  ->  NATHAN / TEST / 2222 / ADMIN
$ws.Range("A4").Value = "NATHAN"
$ws.Range("B4").Value = "TEST"
$ws.Range("C4").NumberFormat = "@"
$ws.Range("C4").Value = "2222"
$ws.Range("D4").Value = "ADMIN"

# Column D width change
$ws.Columns.Item(4).ColumnWidth = 7.710625
